# Commit: "updating all the folders"
# - McCauley rotifers sheet: drop the empty leading "species" column (A),
#   shifting genus/constant/formula left, and lower-case the genus names.
# - Switch the active/selected tab from "McCauley rotifers" to "Bottrell pooled".

$wb = $excel.ActiveWorkbook

# --- McCauley rotifers: remove column A (was always empty aside from the header) ---
$wsRot = $wb.Worksheets.Item("McCauley rotifers")
$wsRot.Activate()
$wsRot.Columns.Item(1).Delete()

# --- Lower-case the rotifer genus names now living in column A ---
$genusNames = @(
    "anuraeopsis",
    "ascomorpha",
    "asplanchna",
    "brachionus",
    "conochilus",
    "collotheca",
    "euchlanis",
    "filinia",
    "gastropus",
    "hexarthra",
    "kellikottia",
    "keratella quadrata",
    "keratella cochlearis",
    "notholca",
    "ploesoma",
    "polyarthra",
    "pompbolyx",
    "synchaeta",
    "testudinella",
    "trichocerca"
)

for ($i = 0; $i -lt $genusNames.Length; $i++) {
    $row = $i + 2
    $wsRot.Cells.Item($row, 1).Value = $genusNames[$i]
}

# Reset the zoom back to 100% and move the selection, matching the saved view state.
$excel.ActiveWindow.Zoom = 100
$wsRot.Range("A6").Select()

# --- Make "Bottrell pooled" the active/selected sheet ---
$wsPooled = $wb.Worksheets.Item("Bottrell pooled")
$wsPooled.Activate()
